$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 388-390 (columns C:F get new values) ---
$ws.Range("C388:F388").Value = 7431065591000
$ws.Range("C389:F389").Value = 7779488501000
$ws.Range("C390:F390").Value = 8162661859000

# --- Append new rows 399-401 ---
# Copy the formatting of the last existing data row (A398) down into the
# new A399:A401 cells so they keep the same style (s="2") as the rest of
# the date column, then set their values.
$ws.Range("A398").Copy()
$ws.Range("A399:A401").PasteSpecial(-4122)

# Row 399
$ws.Range("A399").Value = 44958.45833333334
$ws.Range("B399").Value = "ECONOMICS:ARM2"
$ws.Range("C399:F399").Value = 13551717200000
$ws.Range("G399").Value = 0

# Row 400
$ws.Range("A400").Value = 44986.45833333334
$ws.Range("B400").Value = "ECONOMICS:ARM2"
$ws.Range("C400:F400").Value = 13436944380000
$ws.Range("G400").Value = 0

# Row 401
$ws.Range("A401").Value = 45017.45833333334
$ws.Range("B401").Value = "ECONOMICS:ARM2"
$ws.Range("C401:F401").Value = 14573629300000
$ws.Range("G401").Value = 0
